$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.45386359834331
$ws.Range("C2").Value = 11.76077821269166
$ws.Range("D2").Value = 5.583166910672985
$ws.Range("F2").Value = 29.98433520944133
$ws.Range("G2").Value = 39.03714059998849
$ws.Range("H2").Value = 16.4875912714226
$ws.Range("L2").Value = 10.79867634843321
$ws.Range("M2").Value = 15.78930905851956

$ws.Range("B3").Value = 15.98064158386482
$ws.Range("C3").Value = 11.52026903021662
$ws.Range("D3").Value = 5.596882192471928
$ws.Range("F3").Value = 29.87984104395743
$ws.Range("G3").Value = 38.80074087855212
$ws.Range("H3").Value = 16.51427535156242
$ws.Range("L3").Value = 10.81562366560598
$ws.Range("M3").Value = 15.70976650373311

$ws.Range("B4").Value = 15.68650664293712
$ws.Range("C4").Value = 11.36838073751623
$ws.Range("D4").Value = 5.605684528307354
$ws.Range("F4").Value = 29.82567722347916
$ws.Range("G4").Value = 38.670446715563
$ws.Range("H4").Value = 16.53467756503964
$ws.Range("L4").Value = 10.82760962252298
$ws.Range("M4").Value = 15.6637777759238

$ws.Range("B5").Value = 15.56594060953899
$ws.Range("C5").Value = 11.30547259281961
$ws.Range("D5").Value = 5.609367827779205
$ws.Range("F5").Value = 29.80612979802328
$ws.Range("G5").Value = 38.62112671482807
$ws.Range("H5").Value = 16.54399843873929
$ws.Range("L5").Value = 10.83289129083944
$ws.Range("M5").Value = 15.645768227073

$ws.Range("B6").Value = 15.54588414091189
$ws.Range("C6").Value = 11.294967029222
$ws.Range("D6").Value = 5.609985264880944
$ws.Range("F6").Value = 29.80303670562738
$ws.Range("G6").Value = 38.61316624216537
$ws.Range("H6").Value = 16.54560686030759
$ws.Range("L6").Value = 10.83379230049367
$ws.Range("M6").Value = 15.64282232124333

$ws.Range("B7").Value = 15.68488322997109
$ws.Range("C7").Value = 11.36753637108714
$ws.Range("D7").Value = 5.605733812179035
$ws.Range("F7").Value = 29.82540336550072
$ws.Range("G7").Value = 38.66976623310241
$ws.Range("H7").Value = 16.53479919811871
$ws.Range("L7").Value = 10.82767924444959
$ws.Range("M7").Value = 15.6635319138541

$ws.Range("B8").Value = 16.29155209169523
$ws.Range("C8").Value = 11.67875303041779
$ws.Range("D8").Value = 5.587817110680755
$ws.Range("F8").Value = 29.94623880173568
$ws.Range("G8").Value = 38.95257310314101
$ws.Range("H8").Value = 16.4959559385676
$ws.Range("L8").Value = 10.80419177535594
$ws.Range("M8").Value = 15.76130004134471

$ws.Range("B9").Value = 17.44489829843151
$ws.Range("C9").Value = 12.25366008906707
$ws.Range("D9").Value = 5.555686197037229
$ws.Range("F9").Value = 30.26184293155849
$ws.Range("G9").Value = 39.62288816858943
$ws.Range("H9").Value = 16.45181224026495
$ws.Range("L9").Value = 10.77067365482039
$ws.Range("M9").Value = 15.97496732476785

$ws.Range("B10").Value = 18.26035802224113
$ws.Range("C10").Value = 12.65205892426139
$ws.Range("D10").Value = 5.533882817772722
$ws.Range("F10").Value = 30.54051895642674
$ws.Range("G10").Value = 40.18253480760983
$ws.Range("H10").Value = 16.43908430223103
$ws.Range("L10").Value = 10.75369487920277
$ws.Range("M10").Value = 16.14436539046845

$ws.Range("B11").Value = 18.62260862004869
$ws.Range("C11").Value = 12.82763684220212
$ws.Range("D11").Value = 5.524349589642501
$ws.Range("F11").Value = 30.67714005561361
$ws.Range("G11").Value = 40.45085584525594
$ws.Range("H11").Value = 16.43760326148141
$ws.Range("L11").Value = 10.74763073334037
$ws.Range("M11").Value = 16.22390433186788

$ws.Range("B12").Value = 18.75839999042044
$ws.Range("C12").Value = 12.89327575474932
$ws.Range("D12").Value = 5.52079456213892
$ws.Range("F12").Value = 30.73025933113162
$ws.Range("G12").Value = 40.55435621923001
$ws.Range("H12").Value = 16.43766384788666
$ws.Range("L12").Value = 10.74557289170154
$ws.Range("M12").Value = 16.25436048578077

$ws.Range("B13").Value = 18.72921847299388
$ws.Range("C13").Value = 12.87917751802623
$ws.Range("D13").Value = 5.521557761065705
$ws.Range("F13").Value = 30.71875815019759
$ws.Range("G13").Value = 40.53198272729668
$ws.Range("H13").Value = 16.4376231373851
$ws.Range("L13").Value = 10.74600547916689
$ws.Range("M13").Value = 16.24778656678134

$ws.Range("B14").Value = 18.63380873326276
$ws.Range("C14").Value = 12.83305419585461
$ws.Range("D14").Value = 5.524056015378506
$ws.Range("F14").Value = 30.68148267402452
$ws.Range("G14").Value = 40.45933343168996
$ws.Range("H14").Value = 16.43759578290608
$ws.Range("L14").Value = 10.74745665413078
$ws.Range("M14").Value = 16.22640335072584

$ws.Range("B15").Value = 18.57518340576327
$ws.Range("C15").Value = 12.80469080881758
$ws.Range("D15").Value = 5.525593419109603
$ws.Range("F15").Value = 30.65882956394509
$ws.Range("G15").Value = 40.41507764855761
$ws.Range("H15").Value = 16.43766000141201
$ws.Range("L15").Value = 10.7483765981314
$ws.Range("M15").Value = 16.21334873119672

$ws.Range("B16").Value = 18.23649722583122
$ws.Range("C16").Value = 12.64046760358028
$ws.Range("D16").Value = 5.534513554458877
$ws.Range("F16").Value = 30.53178593940365
$ws.Range("G16").Value = 40.16526887995045
$ws.Range("H16").Value = 16.43926796277119
$ws.Range("L16").Value = 10.75412457088136
$ws.Range("M16").Value = 16.13921552026755

$ws.Range("B17").Value = 18.02639549905459
$ws.Range("C17").Value = 12.53824741284589
$ws.Range("D17").Value = 5.540084151933777
$ws.Range("F17").Value = 30.456349428945
$ws.Range("G17").Value = 40.01548038584705
$ws.Range("H17").Value = 16.44135944024446
$ws.Range("L17").Value = 10.75807575265692
$ws.Range("M17").Value = 16.09435806863325

$ws.Range("B18").Value = 17.90473819793555
$ws.Range("C18").Value = 12.4789226606874
$ws.Range("D18").Value = 5.543324496037727
$ws.Range("F18").Value = 30.41388962792872
$ws.Range("G18").Value = 39.93062406189453
$ws.Range("H18").Value = 16.44296790193485
$ws.Range("L18").Value = 10.76050458707199
$ws.Range("M18").Value = 16.06879180998562

$ws.Range("B19").Value = 17.86341187570444
$ws.Range("C19").Value = 12.45874634331208
$ws.Range("D19").Value = 5.544427864880551
$ws.Range("F19").Value = 30.39967401535877
$ws.Range("G19").Value = 39.9021185281079
$ws.Range("H19").Value = 16.44358208237148
$ws.Range("L19").Value = 10.76135378240243
$ws.Range("M19").Value = 16.06017639977632

$ws.Range("B20").Value = 18.04884626690338
$ws.Range("C20").Value = 12.54918408579343
$ws.Range("D20").Value = 5.539487400032714
$ws.Range("F20").Value = 30.46428382316183
$ws.Range("G20").Value = 40.03129186172427
$ws.Range("H20").Value = 16.44109481701355
$ws.Range("L20").Value = 10.75763897528674
$ws.Range("M20").Value = 16.09910908779401

$ws.Range("B21").Value = 18.66187147384897
$ws.Range("C21").Value = 12.84662503370009
$ws.Range("D21").Value = 5.523320727999361
$ws.Range("F21").Value = 30.69239410323788
$ws.Range("G21").Value = 40.4806215929125
$ws.Range("H21").Value = 16.43758694007387
$ws.Range("L21").Value = 10.74702393683345
$ws.Range("M21").Value = 16.23267514947725

$ws.Range("B22").Value = 19.05439290393107
$ws.Range("C22").Value = 13.03605888369268
$ws.Range("D22").Value = 5.513075270909876
$ws.Range("F22").Value = 30.84952612774135
$ws.Range("G22").Value = 40.78527717336103
$ws.Range("H22").Value = 16.43891700702115
$ws.Range("L22").Value = 10.74147656603601
$ws.Range("M22").Value = 16.32191994348467

$ws.Range("B23").Value = 18.84568102272631
$ws.Range("C23").Value = 12.93541939126701
$ws.Range("D23").Value = 5.518514277236916
$ws.Range("F23").Value = 30.76493690070899
$ws.Range("G23").Value = 40.62169926274913
$ws.Range("H23").Value = 16.43787516749601
$ws.Range("L23").Value = 10.74431015841693
$ws.Range("M23").Value = 16.27411647705866

$ws.Range("B24").Value = 18.03869896507884
$ws.Range("C24").Value = 12.54424134822633
$ws.Range("D24").Value = 5.539757073986434
$ws.Range("F24").Value = 30.46069384761185
$ws.Range("G24").Value = 40.02413956542479
$ws.Range("H24").Value = 16.44121318851053
$ws.Range("L24").Value = 10.75783595252621
$ws.Range("M24").Value = 16.09696045625072

$ws.Range("B25").Value = 17.1378360094302
$ws.Range("C25").Value = 12.10218061465093
$ws.Range("D25").Value = 5.564059810498868
$ws.Range("F25").Value = 30.16814764309439
$ws.Range("G25").Value = 39.42948257932333
$ws.Range("H25").Value = 16.46030544962439
$ws.Range("L25").Value = 10.77839823266728
$ws.Range("M25").Value = 15.91491135200492
